$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recover dropped data: column B ("ID Competição") values were stored without
# the leading "2", so restore them from 69 to 269 for all data rows (2-81).
$ws.Range("B2:B81").Value = 269
